$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date serial values shared by all new rows (Fecha Obtencion / Fecha Publicacion)
$fechaObtencion = 45636.99947916667
$fechaPublicacion = 45545.99947916667

# Populate new rows 24-36 scraped from the PJUD / economicos.cl search

# Row 24
$ws.Range("B24").Value = 'https://www.economicos.cl/remates/clasificados-remates-cod47442289.html'
$ws.Range("C24").Value = $fechaObtencion
$ws.Range("D24").Value = $fechaPublicacion
$ws.Range("F24").Value = 'C-10640-2023'
$ws.Range("G24").Value = '3° juzgado civil de santiago'
$ws.Range("P24").Value = 'vale vista'
$ws.Range("Q24").Value = '10% del mínimo'
$ws.Range("S24").Value = 'Mínimo para iniciar posturas $ 49.736.202'

# Row 25
$ws.Range("B25").Value = 'https://www.economicos.cl/remates/clasificados-remates-cod7478522.html'
$ws.Range("C25").Value = $fechaObtencion
$ws.Range("D25").Value = $fechaPublicacion
$ws.Range("F25").Value = 'N/A'
$ws.Range("G25").Value = 'N/A'
$ws.Range("P25").Value = 'N/A'
$ws.Range("Q25").Value = 'Garantía 10%'
$ws.Range("S25").Value = 'Mínimo $ 37.506.795'

# Row 26
$ws.Range("B26").Value = 'https://www.economicos.cl/remates/clasificados-remates-cod47442024.html'
$ws.Range("C26").Value = $fechaObtencion
$ws.Range("D26").Value = $fechaPublicacion
$ws.Range("F26").Value = 'N/A'
$ws.Range("G26").Value = 'N/A'
$ws.Range("P26").Value = 'vale vista'
$ws.Range("Q26").Value = 'N/A'
$ws.Range("S26").Value = 'N/A'

# Row 27
$ws.Range("B27").Value = 'https://www.economicos.cl/remates/clasificados-remates-cod47444967.html'
$ws.Range("C27").Value = $fechaObtencion
$ws.Range("D27").Value = $fechaPublicacion
$ws.Range("F27").Value = 'N/A'
$ws.Range("G27").Value = 'N/A'
$ws.Range("P27").Value = 'vale vista'
$ws.Range("Q27").Value = 'Garantía 10%'
$ws.Range("S27").Value = 'Mínimo $40.000.000'

# Row 28
$ws.Range("B28").Value = 'https://www.economicos.cl/remates/clasificados-remates-cod47463302.html'
$ws.Range("C28").Value = $fechaObtencion
$ws.Range("D28").Value = $fechaPublicacion
$ws.Range("F28").Value = 'c-1275-2023'
$ws.Range("G28").Value = 'N/A'
$ws.Range("P28").Value = 'N/A'
$ws.Range("Q28").Value = 'N/A'
$ws.Range("S28").Value = 'N/A'

# Row 29
$ws.Range("B29").Value = 'https://www.economicos.cl/remates/clasificados-remates-cod7478514.html'
$ws.Range("C29").Value = $fechaObtencion
$ws.Range("D29").Value = $fechaPublicacion
$ws.Range("F29").Value = 'N/A'
$ws.Range("G29").Value = 'N/A'
$ws.Range("P29").Value = 'Vale Vista'
$ws.Range("Q29").Value = 'caución de un 10%'
$ws.Range("S29").Value = 'Mínimo para la subasta es la suma de $ 138.000.000'

# Row 30
$ws.Range("B30").Value = 'https://www.economicos.cl/remates/clasificados-remates-cod7478519.html'
$ws.Range("C30").Value = $fechaObtencion
$ws.Range("D30").Value = $fechaPublicacion
$ws.Range("F30").Value = 'N/A'
$ws.Range("G30").Value = 'N/A'
$ws.Range("P30").Value = 'N/A'
$ws.Range("Q30").Value = 'Garantía 10%'
$ws.Range("S30").Value = 'Mínimo $ 80.000.000'

# Row 31
$ws.Range("B31").Value = 'https://www.economicos.cl/remates/clasificados-remates-cod47419724.html'
$ws.Range("C31").Value = $fechaObtencion
$ws.Range("D31").Value = $fechaPublicacion
$ws.Range("F31").Value = 'C-3627-2023'
$ws.Range("G31").Value = '1° juzgado de letras de melipilla'
$ws.Range("P31").Value = 'vale vista'
$ws.Range("Q31").Value = '10% del mínimo'
$ws.Range("S31").Value = 'mínimo postura: $ 244.000.000'

# Row 32
$ws.Range("B32").Value = 'https://www.economicos.cl/remates/clasificados-remates-cod47445130.html'
$ws.Range("C32").Value = $fechaObtencion
$ws.Range("D32").Value = $fechaPublicacion
$ws.Range("F32").Value = 'C-4117-2024'
$ws.Range("G32").Value = '2° juzgado civil de santiago'
$ws.Range("P32").Value = 'vale vista'
$ws.Range("Q32").Value = 'N/A'
$ws.Range("S32").Value = 'mínimo para la subasta será la suma de $ 250.838.785'

# Row 33
$ws.Range("B33").Value = 'https://www.economicos.cl/remates/clasificados-remates-cod7477212.html'
$ws.Range("C33").Value = $fechaObtencion
$ws.Range("D33").Value = $fechaPublicacion
$ws.Range("F33").Value = 'N/A'
$ws.Range("G33").Value = 'N/A'
$ws.Range("P33").Value = 'N/A'
$ws.Range("Q33").Value = 'Garantía 10%'
$ws.Range("S33").Value = 'Mínimo $ 45.000.000'

# Row 34
$ws.Range("B34").Value = 'https://www.economicos.cl/remates/clasificados-remates-cod47442147.html'
$ws.Range("C34").Value = $fechaObtencion
$ws.Range("D34").Value = $fechaPublicacion
$ws.Range("F34").Value = 'N/A'
$ws.Range("G34").Value = 'N/A'
$ws.Range("P34").Value = 'N/A'
$ws.Range("Q34").Value = 'Garantía 10%'
$ws.Range("S34").Value = 'Mínimo UF 3.567'

# Row 35
$ws.Range("B35").Value = 'https://www.economicos.cl/remates/clasificados-remates-cod47445033.html'
$ws.Range("C35").Value = $fechaObtencion
$ws.Range("D35").Value = $fechaPublicacion
$ws.Range("F35").Value = 'C-26213-2019'
$ws.Range("G35").Value = '6° juzgado civil de santiago'
$ws.Range("P35").Value = 'vale vista'
$ws.Range("Q35").Value = 'N/A'
$ws.Range("S35").Value = 'mínimo para iniciar las posturas: $ 91.901.468'

# Row 36
$ws.Range("B36").Value = 'https://www.economicos.cl/remates/clasificados-remates-cod47458046.html'
$ws.Range("C36").Value = $fechaObtencion
$ws.Range("D36").Value = $fechaPublicacion
$ws.Range("F36").Value = 'N/A'
$ws.Range("G36").Value = 'N/A'
$ws.Range("P36").Value = 'vale vista'
$ws.Range("Q36").Value = 'interesados equivalente al 10%'
$ws.Range("S36").Value = 'Mínimo: $70.308.933'

# Apply the same date number format (style used by C5:D23) to the new date cells
$ws.Range("C23:D23").Copy()
$ws.Range("C24:D36").PasteSpecial(-4122)  # xlPasteFormats

# Mark the range B5:S37 as containing text that should not be flagged as "number stored as text"
try {
    $ws.Range("B5:S37").Errors.Item(9).Ignore = $true
} catch {}

# Ensure the worksheet dimension / used range extends through row 37 to match the source workbook
$ws.Range("S37").Font.Bold = $false
